# "aggiornamento fino al 26/03" - append new daily rows (21/04-26/04/2021
# serials 44308-44312) below the existing data on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 233

# Carry the date-column formatting (style/number format) from the last
# existing row down onto the new rows before writing values.
$ws.Range("A$lastRow").Copy() | Out-Null
$ws.Range("A234:A238").PasteSpecial(-4122) | Out-Null

$newRows = @(
    @(44308, 1, 4, 45.99816007359706),
    @(44309, 1, 4, 45.99816007359706),
    @(44310, 2, 6, 68.99724011039559),
    @(44311, 2, 7, 80.49678012879485),
    @(44312, 0, 7, 80.49678012879485)
)

$row = $lastRow + 1
foreach ($vals in $newRows) {
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $row++
}
